# Replace the crew-member name "Leo Croufer" with "Justin Ferrandez" in the
# signature table (bloc B of the form).
#
# The two words are currently split across two separate <w:r> runs ("Leo "
# and "Croufer", the latter wrapped in <w:proofErr> spell-check markers
# because "Croufer" was flagged as a misspelling). The target content has
# a single run with the literal text "Justin Ferrandez" and no proofErr
# wrapper, so both the text and the run layout need to change.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "Leo Croufer"
$find.Execute() | Out-Null

if ($find.Found) {
    $target = $find.Parent
    $start = $target.Start
    $newText = "Justin Ferrandez"
    $wasBold = $target.Bold

    # A plain Range.Text assignment here would keep the formatting identical
    # to the immediately-preceding run ("...Patron d'embarcation ; ") and the
    # engine then silently coalesces the two into a single run. Toggling a
    # direct-formatting flag across the edit keeps this run distinct, then
    # we restore the flag to its original value afterwards.
    $target.Bold = 1
    $target.Text = $newText

    $after = $d.Range($start, $start + $newText.Length)
    $after.Bold = $wasBold
}
